$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace the Min-column formulas (=D/5) with hard values equal to Max,
# per the commit "Tentativas de corrigir o problema dos adopters negativos".
$ws.Range("C2").Value = 0.01
$ws.Range("C3").Value = 100
$ws.Range("C4").Value = 0.02

# aTotalPopulation: bump Max to 1000000 and mirror it into Min (no formula).
$ws.Range("D5").Value = 1000000
$ws.Range("C5").Value = 1000000

# Drop the siniPotentialAdopters / "Consumidores Potenciais" row entirely.
$ws.Range("A6:E6").ClearContents()

# Restore the selection to B4, matching the saved view state.
$ws.Range("B4").Select() | Out-Null
